$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet right after "raw data"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "the graph data"

# Header row
$ws2.Range("B1").Value2 = "JD - 2457651.0[day]"
$ws2.Range("C1").Value2 = "V[mag]"
$ws2.Range("D1").Value2 = "error_V[mag]"

# Copy the header style (bold, centered, bordered) from the "raw data" sheet
$ws1.Range("B1").Copy()
$ws2.Range("B1:D1").PasteSpecial(-4122)
$ws2.Range("A2:A50").PasteSpecial(-4122)

$data = @(
  @(0.0, 0.6650000000372529, 17.022, 0.0),
  @(1.0, 1.799000000115484, 14.672, 0.04),
  @(2.0, 1.799000000115484, 14.672, 0.04),
  @(3.0, 2.165000000037253, 14.789, 0.064),
  @(4.0, 2.231000000145286, 14.957, 0.061),
  @(5.0, 2.606999999843538, 15.106, 0.012),
  @(6.0, 2.608000000007451, 15.102, 0.008),
  @(7.0, 2.62900000018999, 15.202, 0.03),
  @(8.0, 2.62900000018999, 15.202, 0.03),
  @(9.0, 2.788000000175089, 15.392, 0.05),
  @(10.0, 2.788000000175089, 15.392, 0.05),
  @(11.0, 2.820999999996275, 15.274, 0.012),
  @(12.0, 2.822000000160187, 15.247, 0.014),
  @(13.0, 2.896000000182539, 15.502, 0.06),
  @(14.0, 2.896000000182539, 15.502, 0.06),
  @(15.0, 3.234999999869615, 15.652, 0.075),
  @(16.0, 3.31900000013411, 15.756, 0.038),
  @(17.0, 3.319999999832362, 15.949, 0.034),
  @(18.0, 3.572999999858439, 16.023, 0.016),
  @(19.0, 3.574000000022352, 16.044, 0.018),
  @(20.0, 3.666999999899417, 16.059, 0.021),
  @(21.0, 3.66800000006333, 16.085, 0.02),
  @(22.0, 3.697999999858439, 16.115, 0.023),
  @(23.0, 3.699000000022352, 16.123, 0.018),
  @(24.0, 3.873999999836087, 16.208, 0.021),
  @(25.0, 3.875, 16.187, 0.022),
  @(26.0, 3.88799999980256, 16.282, 0.015),
  @(27.0, 3.888999999966472, 16.289, 0.014),
  @(28.0, 4.320999999996275, 16.235, 0.024),
  @(29.0, 4.322000000160187, 16.345, 0.031),
  @(30.0, 4.688999999780208, 16.372, 0.018),
  @(31.0, 4.689999999944121, 16.323, 0.02),
  @(32.0, 4.728999999817461, 16.319, 0.018),
  @(33.0, 4.729999999981374, 16.29, 0.018),
  @(34.0, 5.313000000081956, 16.156, 0.047),
  @(35.0, 5.313999999780208, 16.265, 0.044),
  @(36.0, 5.458999999798834, 16.111, 0.04),
  @(37.0, 5.459999999962747, 16.131, 0.041),
  @(38.0, 5.518000000156462, 16.141, 0.014),
  @(39.0, 5.518999999854714, 16.185, 0.012),
  @(40.0, 5.546000000089407, 16.197, 0.044),
  @(41.0, 5.546999999787658, 16.224, 0.039),
  @(42.0, 6.416000000201166, 15.949, 0.009),
  @(43.0, 6.41800000006333, 15.935, 0.009),
  @(44.0, 6.504999999888241, 15.881, 0.017),
  @(45.0, 6.506000000052154, 15.925, 0.017),
  @(46.0, 6.580000000074506, 15.965, 0.009),
  @(47.0, 6.584999999962747, 15.837, 0.028),
  @(48.0, 6.58600000012666, 15.906, 0.028)
)

foreach ($row in $data) {
  $r = [int]$row[0] + 2
  $ws2.Cells.Item($r, 1).Value2 = $row[0]
  $ws2.Cells.Item($r, 2).Value2 = $row[1]
  $ws2.Cells.Item($r, 3).Value2 = $row[2]
  $ws2.Cells.Item($r, 4).Value2 = $row[3]
}

# Column widths for C and D (best fit)
$ws2.Columns.Item(3).ColumnWidth = 7.19921875
$ws2.Columns.Item(4).ColumnWidth = 12.5

# Selections
$ws1.Range("C44").Select()
$ws2.Range("H10").Select()
$ws2.Activate()

Write-Output "done"
